$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-26 Thursday", "2024-12-27 Friday"),
    @("326÷2=", "553÷7="),
    @("469÷9=", "888÷9="),
    @("777÷6=", "893÷6="),
    @("439÷6=", "626÷6="),
    @("692÷5=", "445÷3="),
    @("146÷8=", "226÷7="),
    @("492÷6=", "748÷2="),
    @("835÷4=", "716÷9="),
    @("109÷2=", "160÷8="),
    @("692÷2=", "978÷8="),
    @("188÷9=", "753÷8="),
    @("169÷9=", "354÷2="),
    @("359÷5=", "892÷3="),
    @("329÷6=", "588÷8="),
    @("990÷3=", "791÷5="),
    @("220÷4=", "699÷4="),
    @("935÷2=", "778÷9="),
    @("819÷8=", "804÷4="),
    @("479÷4=", "172÷8="),
    @("577÷5=", "767÷6="),
    @("434÷8=", "453÷7="),
    @("322÷6=", "949÷8="),
    @("869÷3=", "850÷6="),
    @("269÷7=", "643÷2="),
    @("141÷4=", "314÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
